$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "source" column (F) to make room
# for the new "file_size" field, shifting source/link/notes one column right.
$ws.Range("F1").EntireColumn.Insert()

# Populate the new header cell.
$ws.Range("F1").Value = "file_size"

# Match the column width used by similarly-typed text columns (e.g. E).
$ws.Range("F1").EntireColumn.ColumnWidth = 9.29

# Update the active selection to match the post-edit workbook state.
$ws.Range("F2").Select()
